$d = $word.ActiveDocument

# Locate the existing "NSC = Neural Statistics Compiler" abbreviation
# list item in the Abbreviations section, and append two new list
# items after it, matching the ListParagraph / numId=2 / ilvl=0 style.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "NSC = Neural Statistics Compiler") {
        $anchor = $p
        break
    }
}

# First new paragraph: CTB abbreviation (single run)
$anchor.Range.InsertParagraphAfter()
$ctbXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>CTB = Cell Titer Blue. Same assay as Alamar Blue</w:t></w:r></w:p>'
$ctbPara = $anchor.Next()
[void]$ctbPara.Range.InsertXML($ctbXml)

# Second new paragraph: AB abbreviation (split across two runs, as in
# the source document)
$ctbPara.Range.InsertParagraphAfter()
$abXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>AB  = Alamar Blue. Same assay as C</w:t></w:r><w:r><w:t>ell Titer Blue</w:t></w:r></w:p>'
$abPara = $ctbPara.Next()
[void]$abPara.Range.InsertXML($abXml)

Write-Host "Abbreviations updated."
